$d = $word.ActiveDocument

# Locate the table that holds the "id / title / author / ..." header row
$t = $d.Tables.Item(1)
$tblStart = $t.Range.Start

# Delete the table; the paragraph mark that used to follow it now sits at $tblStart
$null = $t.Delete()

# Replace that paragraph (its text run plus its end-of-paragraph mark) with a
# single paragraph containing one bold run per field, tab-separated. A second,
# empty paragraph is appended in the same InsertXML call to restore the extra
# blank paragraph mark that gets consumed by the replace, so the two blank
# paragraphs that originally followed the table are preserved.
$ins = $d.Range($tblStart, $tblStart + 1)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> id  </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/><w:t xml:space="preserve"> title  </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/><w:t xml:space="preserve"> author  </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/><w:t xml:space="preserve"> published_date  </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/><w:t xml:space="preserve"> isbn  </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/><w:t xml:space="preserve"> summary  </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/><w:t xml:space="preserve"> image  </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/><w:t xml:space="preserve"> images1  </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/><w:t xml:space="preserve"></w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$null = $ins.InsertXML($xml)
